$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.243.51'
$ws.Range("E2").Value = '  +0.37%  '

$ws.Range("D3").Value = '1.683.58'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.72'
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5262'
$ws.Range("E6").Value = '  +3.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2698'
$ws.Range("E8").Value = '  +1.98%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06431'
$ws.Range("E9").Value = '  +1.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.05'
$ws.Range("E10").Value = '  +1.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07452'
$ws.Range("E11").Value = '  +0.98%  '

$ws.Range("D12").Value = '1.686.91'
$ws.Range("E12").Value = '  +0.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.539'
$ws.Range("E13").Value = '  +0.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5865'
$ws.Range("E14").Value = '  +1.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008604'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.67'
$ws.Range("E16").Value = '  -0.30%  '

$ws.Range("D17").Value = '26.265.55'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.968'
$ws.Range("E18").Value = '  -0.46%  '

$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.83'
$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.90'
$ws.Range("E21").Value = '  +2.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.238'
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  -0.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.34'
$ws.Range("E24").Value = '  +1.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1251'
$ws.Range("E25").Value = '  +7.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.668'
$ws.Range("E26").Value = '  +1.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.84'
$ws.Range("E27").Value = '  +0.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06650'
$ws.Range("E28").Value = '  +15.36%  '

$ws.Range("E29").Value = '  +1.43%  '

$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("E31").Value = '  +2.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.555'
$ws.Range("E32").Value = '  +1.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.671'
$ws.Range("E33").Value = '  +0.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.023'
$ws.Range("E34").Value = '  +2.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6209'
$ws.Range("E35").Value = '  +3.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.371'
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  +2.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.284'
$ws.Range("E38").Value = '  +6.27%  '

$ws.Range("D39").Value = '1.100.60'
$ws.Range("E39").Value = '  +0.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01607'
$ws.Range("E40").Value = '  +0.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8718'
$ws.Range("E41").Value = '  +1.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.012'
$ws.Range("E42").Value = '  +0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.20'
$ws.Range("E43").Value = '  +1.96%  '

$ws.Range("D44").Value = '1.827.63'
$ws.Range("E44").Value = '  +0.45%  '

$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.79'
$ws.Range("E46").Value = '  +1.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.180'
$ws.Range("E47").Value = '  +1.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05250'
$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4282'
$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.017'
$ws.Range("E51").Value = '  +3.38%  '
